# Update stats for 2025-08 (row 21, month=45870)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = 6223
$ws.Range("C21").Value = 981
$ws.Range("D21").Value = 5584721
$ws.Range("E21").Value = 897.4322673951471
$ws.Range("F21").Value = 8.019441069258804
$ws.Range("G21").Value = 3.481012658227844
$ws.Range("H21").Value = 27.47675629974353
